# Rename header labels to match the "Fixed update to excel issue" commit
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new "PO Forecast" sheet after "Monthly Trend"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$hdr = $wsForecast.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# Data rows
$forecastData = @(
    @(45312.99999999999, 76, -5.825370827157596, 153.8978512955196),
    @(45319.99999999999, 77, -0.05423358162458949, 152.2165448565046),
    @(45361.99999999999, 85, 2.576470551037624, 164.8781268990627),
    @(45389.99999999999, 91, 14.52493389531118, 166.1462253415498),
    @(45403.99999999999, 94, 19.73947797567184, 171.3878286617272),
    @(45417.99999999999, 97, 15.28190853242381, 173.3824606870184),
    @(45424.99999999999, 98, 18.78975855923144, 175.5743032791049),
    @(45431.99999999999, 100, 14.71386830094464, 174.2002910175407),
    @(45438.99999999999, 101, 22.3198871630005, 175.6456431451187),
    @(45445.99999999999, 102, 23.87968813020984, 176.2557175143116),
    @(45452.99999999999, 104, 25.66308661182341, 180.8632994418697),
    @(45459.99999999999, 105, 29.60368719736113, 190.6707604886587),
    @(45473.99999999999, 108, 31.78908826020726, 183.1267607866607),
    @(45480.99999999999, 109, 31.27201804964945, 196.586082129481),
    @(45487.99999999999, 111, 31.25972424362431, 185.6245627658288),
    @(45494.99999999999, 112, 30.5804950400725, 193.4922847514086),
    @(45501.99999999999, 114, 31.48818979976226, 195.3252407623887),
    @(45508.99999999999, 115, 40.28004695052734, 193.3765953458056),
    @(45515.99999999999, 116, 39.23365500536912, 196.870342212232),
    @(45529.99999999999, 119, 39.94112507269492, 196.6713752373882),
    @(45536.99999999999, 121, 44.88924558723504, 199.8699946854641),
    @(45557.99999999999, 125, 43.6283609559161, 206.1061516840991),
    @(45571.99999999999, 128, 53.38786224559009, 206.0336405606417),
    @(45599.99999999999, 133, 54.849388764006, 214.5419317996872),
    @(45606.99999999999, 135, 54.11033806228011, 210.3321324042048),
    @(45613.99999999999, 136, 58.95759509525849, 218.3345603410376),
    @(45620.99999999999, 138, 57.39007002693189, 218.1564439297005),
    @(45627.99999999999, 139, 65.09637594104956, 219.5010676455826),
    @(45634.99999999999, 140, 60.33783421355759, 211.1264760207514),
    @(45641.99999999999, 142, 71.50763282588528, 218.1045976284379),
    @(45648.99999999999, 143, 63.1591030816551, 220.3827117880992),
    @(45655.99999999999, 145, 64.62692772286204, 221.9357337582354)
)

$r = 2
foreach ($row in $forecastData) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r++
}

$wsForecast.Range("A1").Select()
Write-Host "PO Forecast sheet populated with $($forecastData.Count) data rows"
